$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: add the TARGET commit hash in column C (reuses column C's
# default "no fill" style, then gets the 0.00E+00 number format applied
# to match the author's edit).
$ws.Range("C9").Value2 = "545e905dcec09f579fb05c55beab5bf79537641a"
$ws.Range("C9").NumberFormat = "0.00E+00"

# Row 9: flip STATUS from "open" to "closed" and flag it with a solid
# red fill (new fill color, distinct from the existing "closed" rows).
$ws.Range("B9").Value2 = "closed"
$ws.Range("B9").Interior.Color = 255

# Move the active selection down to B10, matching the saved view state.
$ws.Range("B10").Select() | Out-Null
